# Applies the "Latest data" refresh to the particelle_non_trovate sheet.
# The data table (columns B:C, rows 2..73) is replaced with the values
# that used to live 3 rows further down (rows 5..76) in the previous
# version of the sheet, and the now-unused trailing rows 74..76 are
# removed so the sheet shrinks from A1:C76 to A1:C73.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for column B (codice_particella) and C (codice_comune_catastale)
# for rows 2..73 of the resulting sheet.
$newB = @(
    "2181/2","9769/2","9783/231","13823","14724/2","8974/1","15380","466",
    ".315","1900/4","754","755","825/63","1272/3","765","94/6",
    "2681/1","789","384/1","1117/2","1230/100","1230/115","1230/85","1230/86",
    "1230/87","1230/88","1303/1","1303/2","1309","1330","1334","1346",
    "1369/1","194/4","254/2","337/5","393/1","393/2","393/3","465",
    "614","1038/2","420/101","420/102","420/106","420/107","420/109","420/110",
    "420/80","420/92","420/93","420/94","420/95","420/96","420/97","454",
    "53","705/11","756","798/3","4523","1911/5","1912/5","1912/5",
    "3597/16","3597/18","3597/22","2331/38","2068/43","2822/12","2822/16","2020/14"
)

$newC = @(
    57,44,44,154,154,154,154,213,
    282,317,97,97,97,390,404,251,
    442,443,193,193,193,193,193,193,
    193,193,193,193,193,193,193,193,
    193,193,193,193,193,193,193,193,
    193,215,215,215,215,215,215,215,
    215,215,215,215,215,215,215,215,
    215,215,215,215,404,404,404,404,
    9,9,9,258,310,310,310,310
)

$startRow = 2
$endRow = 73

if ($newB.Count -ne ($endRow - $startRow + 1)) {
    throw "newB count mismatch: $($newB.Count)"
}
if ($newC.Count -ne ($endRow - $startRow + 1)) {
    throw "newC count mismatch: $($newC.Count)"
}

# Force column B to be treated as text so purely-numeric-looking particle
# codes (e.g. "13823", ".315") are not silently converted to numbers.
$colB = $ws.Range("B$startRow`:B$endRow")
$colB.NumberFormat = "@"

for ($i = 0; $i -lt $newB.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 2).Value = $newB[$i]
    $ws.Cells.Item($r, 3).Value = $newC[$i]
}

# Restore the default (unstyled) look for column B now that the text
# values are safely stored, so no extra number-format style lingers.
$colB.Style = "Normal"

# Drop the three now-obsolete trailing rows (old rows 74-76), shrinking
# the sheet's used range from A1:C76 to A1:C73.
$ws.Rows("74:76").Delete() | Out-Null
